$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1761.5714
$ws.Range("I4").Value = 2070.4
$ws.Range("J4").Value = 989.5
$ws.Range("K4").Value = 2070.4
$ws.Range("L4").Value = 989.5
$ws.Range("M4").Value = -1956.4
$ws.Range("N4").Value = -1217.5
$ws.Range("H116").Value = 11115140
$ws.Range("I116").Value = 18522238
$ws.Range("K116").Value = 18522238
$ws.Range("M116").Value = -18518796
$ws.Range("H125").Value = 5041.2
$ws.Range("I125").Value = 6337
$ws.Range("J125").Value = 3097.5
$ws.Range("K125").Value = 57033
$ws.Range("L125").Value = 27877.5
$ws.Range("M125").Value = -54573
$ws.Range("N125").Value = -32797.5
$ws.Range("H132").Value = 2275794.8
$ws.Range("I132").Value = 3172.9736
$ws.Range("K132").Value = 9518.9208
$ws.Range("M132").Value = -6988.9208
$ws.Range("H137").Value = 11437.318
$ws.Range("I137").Value = 17439
$ws.Range("K137").Value = 52317
$ws.Range("M137").Value = -49767

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4969.8027
$ws.Range("J32").Value = 8399.75
$ws.Range("L32").Value = 8399.75
$ws.Range("N32").Value = -8973.75
$ws.Range("H45").Value = 207097.5
$ws.Range("I45").Value = 339312.5
$ws.Range("K45").Value = 339312.5
$ws.Range("M45").Value = -338935.5
$ws.Range("H60").Value = 5032.1113
$ws.Range("I60").Value = 4411.125
$ws.Range("K60").Value = 4411.125
$ws.Range("M60").Value = -3678.125
$ws.Range("H61").Value = 17420.416
$ws.Range("I61").Value = 23631.125
$ws.Range("K61").Value = 23631.125
$ws.Range("M61").Value = -23419.125
$ws.Range("H74").Value = 6934.0527
$ws.Range("I74").Value = 15192.429
$ws.Range("J74").Value = 2116.6667
$ws.Range("K74").Value = 15192.429
$ws.Range("L74").Value = 2116.6667
$ws.Range("M74").Value = -14318.429
$ws.Range("N74").Value = -3864.6667
$ws.Range("H77").Value = 6934.0527
$ws.Range("I77").Value = 15192.429
$ws.Range("J77").Value = 2116.6667
$ws.Range("K77").Value = 75962.145
$ws.Range("L77").Value = 10583.3335
$ws.Range("M77").Value = -71594.145
$ws.Range("N77").Value = -19319.3335
$ws.Range("H122").Value = 1770463.1
$ws.Range("I122").Value = 5221.077
$ws.Range("K122").Value = 15663.231
$ws.Range("M122").Value = -13213.231
$ws.Range("H132").Value = 2188.9788
$ws.Range("I132").Value = 1566.0605
$ws.Range("J132").Value = 3657.2856
$ws.Range("K132").Value = 4698.181500000001
$ws.Range("L132").Value = 10971.8568
$ws.Range("M132").Value = -2168.181500000001
$ws.Range("N132").Value = -16031.8568
$ws.Range("H136").Value = 17420.416
$ws.Range("I136").Value = 23631.125
$ws.Range("K136").Value = 70893.375
$ws.Range("M136").Value = -68343.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3471.625
$ws.Range("I20").Value = 1619.4
$ws.Range("K20").Value = 1619.4
$ws.Range("M20").Value = -1372.4
$ws.Range("H96").Value = 7870.2
$ws.Range("I96").Value = 7870.2
$ws.Range("K96").Value = 7870.2
$ws.Range("M96").Value = -5124.2
$ws.Range("H134").Value = 9471.842000000001
$ws.Range("I134").Value = 11844.462
$ws.Range("J134").Value = 4331.1665
$ws.Range("K134").Value = 35533.386
$ws.Range("L134").Value = 12993.4995
$ws.Range("M134").Value = -32998.386
$ws.Range("N134").Value = -18063.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5663.3335
$ws.Range("I31").Value = 5903.3706
$ws.Range("J31").Value = 4583.1665
$ws.Range("K31").Value = 5903.3706
$ws.Range("L31").Value = 4583.1665
$ws.Range("M31").Value = -5608.3706
$ws.Range("N31").Value = -5173.1665
$ws.Range("H34").Value = 5663.3335
$ws.Range("I34").Value = 5903.3706
$ws.Range("J34").Value = 4583.1665
$ws.Range("K34").Value = 5903.3706
$ws.Range("L34").Value = 4583.1665
$ws.Range("M34").Value = -5701.3706
$ws.Range("N34").Value = -4987.1665
$ws.Range("H58").Value = 2806.372
$ws.Range("I58").Value = 2823.7778
$ws.Range("J58").Value = 2777
$ws.Range("K58").Value = 2823.7778
$ws.Range("L58").Value = 2777
$ws.Range("M58").Value = -2620.7778
$ws.Range("N58").Value = -3183
$ws.Range("H122").Value = 14428.777
$ws.Range("I122").Value = 16119.875
$ws.Range("K122").Value = 48359.625
$ws.Range("M122").Value = -45909.625
$ws.Range("H132").Value = 2032
$ws.Range("I132").Value = 2032
$ws.Range("K132").Value = 6096
$ws.Range("M132").Value = -3566
$ws.Range("H136").Value = 2806.372
$ws.Range("I136").Value = 2823.7778
$ws.Range("J136").Value = 2777
$ws.Range("K136").Value = 8471.3334
$ws.Range("L136").Value = 8331
$ws.Range("M136").Value = -5921.3334
$ws.Range("N136").Value = -13431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 501352.56
$ws.Range("J5").Value = 835215.25
$ws.Range("L5").Value = 2505645.75
$ws.Range("N5").Value = -2505869.75
$ws.Range("J62").Value = 1000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4372
$ws.Range("J65").Value = 1000
$ws.Range("L65").Value = 9000
$ws.Range("N65").Value = -15864
$ws.Range("H119").Value = 707
$ws.Range("I119").Value = 707
$ws.Range("K119").Value = 2121
$ws.Range("M119").Value = 2717
$ws.Range("H122").Value = 4622.925
$ws.Range("I122").Value = 1488.3334
$ws.Range("J122").Value = 5176.0884
$ws.Range("K122").Value = 13395.0006
$ws.Range("L122").Value = 46584.7956
$ws.Range("M122").Value = -10945.0006
$ws.Range("N122").Value = -51484.7956
$ws.Range("H135").Value = 501352.56
$ws.Range("J135").Value = 835215.25
$ws.Range("L135").Value = 7516937.25
$ws.Range("N135").Value = -7522007.25
$ws.Range("H137").Value = 2654.1
$ws.Range("I137").Value = 2378.2856
$ws.Range("J137").Value = 3297.6667
$ws.Range("K137").Value = 7134.8568
$ws.Range("L137").Value = 9893.000100000001
$ws.Range("M137").Value = -2034.8568
$ws.Range("N137").Value = -20093.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 15166.333
$ws.Range("I97").Value = 17785.285
$ws.Range("J97").Value = 6000
$ws.Range("K97").Value = 17785.285
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -17289.285
$ws.Range("N97").Value = -6992
$ws.Range("H122").Value = 8566.799999999999
$ws.Range("I122").Value = 5572.609
$ws.Range("K122").Value = 16717.827
$ws.Range("M122").Value = -14267.827
$ws.Range("H132").Value = 10291.182
$ws.Range("I132").Value = 12401.125
$ws.Range("J132").Value = 4664.6665
$ws.Range("K132").Value = 37203.375
$ws.Range("L132").Value = 13993.9995
$ws.Range("M132").Value = -34673.375
$ws.Range("N132").Value = -19053.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 24326.6
$ws.Range("I7").Value = 27266.705
$ws.Range("J7").Value = 7666
$ws.Range("K7").Value = 27266.705
$ws.Range("L7").Value = 7666
$ws.Range("M7").Value = -27154.705
$ws.Range("N7").Value = -7890
$ws.Range("H16").Value = 3352.125
$ws.Range("I16").Value = 2831
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 2831
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = -2661
$ws.Range("N16").Value = -7340
$ws.Range("H40").Value = 34646.35
$ws.Range("I40").Value = 44090.453
$ws.Range("K40").Value = 44090.453
$ws.Range("M40").Value = -43954.453
$ws.Range("H55").Value = 1809.7
$ws.Range("I55").Value = 311.6
$ws.Range("K55").Value = 311.6
$ws.Range("M55").Value = -138.6
$ws.Range("H96").Value = 26666.666
$ws.Range("J96").Value = 26666.666
$ws.Range("L96").Value = 26666.666
$ws.Range("N96").Value = -32158.666
$ws.Range("H122").Value = 6525.625
$ws.Range("J122").Value = 4796.2
$ws.Range("L122").Value = 14388.6
$ws.Range("N122").Value = -19288.6
$ws.Range("H126").Value = 24326.6
$ws.Range("I126").Value = 27266.705
$ws.Range("J126").Value = 7666
$ws.Range("K126").Value = 81800.11500000001
$ws.Range("L126").Value = 22998
$ws.Range("M126").Value = -79330.11500000001
$ws.Range("N126").Value = -27938
$ws.Range("H132").Value = 385629.44
$ws.Range("I132").Value = 599074.4
$ws.Range("K132").Value = 1797223.2
$ws.Range("M132").Value = -1794693.2
$ws.Range("H136").Value = 5504.357
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 22104.088
$ws.Range("I126").Value = 28541.295
$ws.Range("J126").Value = 3865.3333
$ws.Range("K126").Value = 85623.88499999999
$ws.Range("L126").Value = 11595.9999
$ws.Range("M126").Value = -83153.88499999999
$ws.Range("N126").Value = -16535.9999
$ws.Range("H132").Value = 3113.3555
$ws.Range("I132").Value = 2962.0857
$ws.Range("K132").Value = 8886.257100000001
$ws.Range("M132").Value = -6356.257100000001
$ws.Range("H136").Value = 2601053.2
$ws.Range("I136").Value = 3101263.8
$ws.Range("K136").Value = 9303791.399999999
$ws.Range("M136").Value = -9301241.399999999
